$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.343.10'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '3.759.53'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").Value = '3.756.10'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.97%  '
$ws.Range("E12").Value = '  -2.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.13'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '4.389.15'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '3.764.46'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '69.375.49'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("E18").Value = '  -2.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '499.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.47%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.350'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("E39").Value = '  +4.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '463.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.61%  '
$ws.Range("E42").Value = '  -5.14%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.49%  '
$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = '2.951.12'
$ws.Range("E46").Value = '  -3.96%  '
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.31%  '
$ws.Range("E51").Value = '  -1.41%  '
